$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Ref, $Val)
    $c = $ws.Range($Ref)
    $c.NumberFormat = "@"
    $c.Value = $Val
    $c.Style = "Normal"
}

Set-TextValue "D2" "27.938.00"
$ws.Range("E2").Value = "  +0.79%  "
Set-TextValue "D3" "1.764.55"
$ws.Range("E3").Value = "  -0.41%  "
Set-TextValue "D4" "1.001"
$ws.Range("E4").Value = "  -0.13%  "
Set-TextValue "D5" "328.60"
$ws.Range("E5").Value = "  +0.51%  "
Set-TextValue "D6" "1.001"
$ws.Range("E6").Value = "  -0.08%  "
Set-TextValue "D7" "0.4683"
$ws.Range("E7").Value = "  +2.03%  "
Set-TextValue "D8" "0.3522"
$ws.Range("E8").Value = "  -1.60%  "
Set-TextValue "D9" "43.66"
$ws.Range("E9").Value = "  +4.47%  "
Set-TextValue "D10" "0.07372"
Set-TextValue "D11" "1.080"
$ws.Range("E11").Value = "  -1.63%  "
$ws.Range("E12").Value = "  -0.13%  "
Set-TextValue "D13" "20.58"
$ws.Range("E13").Value = "  -0.93%  "
Set-TextValue "D14" "5.992"
$ws.Range("E14").Value = "  -0.67%  "
Set-TextValue "D15" "7.169"
$ws.Range("E15").Value = "  -0.55%  "
Set-TextValue "D16" "1.764.00"
$ws.Range("E16").Value = "  -0.44%  "
Set-TextValue "D17" "92.19"
$ws.Range("E17").Value = "  -1.50%  "
Set-TextValue "D18" "0.00001053"
$ws.Range("E18").Value = "  -0.44%  "
Set-TextValue "D19" "0.06419"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("E20").Value = "  -0.10%  "
Set-TextValue "D21" "16.89"
$ws.Range("E21").Value = "  -0.83%  "
Set-TextValue "D22" "5.769"
$ws.Range("E22").Value = "  -0.34%  "
Set-TextValue "D23" "27.963.45"
$ws.Range("E23").Value = "  +0.64%  "
Set-TextValue "D24" "11.12"
$ws.Range("E24").Value = "  -1.48%  "
Set-TextValue "D25" "2.153"
$ws.Range("E25").Value = "  +3.18%  "
Set-TextValue "D26" "162.46"
$ws.Range("E26").Value = "  -1.10%  "
Set-TextValue "D27" "19.98"
$ws.Range("E27").Value = "  -1.24%  "
Set-TextValue "D28" "1.966.25"
$ws.Range("E28").Value = "  -0.58%  "
Set-TextValue "D29" "2.171"
$ws.Range("E29").Value = "  -0.13%  "
Set-TextValue "D30" "122.77"
$ws.Range("E30").Value = "  -2.43%  "
Set-TextValue "D31" "1.069"
$ws.Range("E31").Value = "  -1.73%  "
Set-TextValue "D32" "0.09288"
$ws.Range("E32").Value = "  +0.80%  "
Set-TextValue "D33" "3.648"
$ws.Range("E33").Value = "  -0.66%  "
Set-TextValue "D34" "5.539"
$ws.Range("E34").Value = "  +0.39%  "
Set-TextValue "D35" "11.66"
$ws.Range("E35").Value = "  -1.37%  "
Set-TextValue "D36" "0.02262"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("E38").Value = "  -1.08%  "
Set-TextValue "D39" "4.901"
$ws.Range("E39").Value = "  -0.75%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D40" "1.187"
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D41" "0.6130"
$ws.Range("E41").Value = "  -2.64%  "
Set-TextValue "D42" "1.419"
$ws.Range("E42").Value = "  +2.00%  "
Set-TextValue "D43" "7.742"
$ws.Range("E43").Value = "  -0.30%  "
Set-TextValue "D44" "13.17"
$ws.Range("E44").Value = "  -0.60%  "
Set-TextValue "D45" "3.738"
$ws.Range("E45").Value = "  -0.08%  "
Set-TextValue "D46" "0.5781"
$ws.Range("E46").Value = "  -1.77%  "
Set-TextValue "D47" "123.32"
$ws.Range("E47").Value = "  +0.80%  "
Set-TextValue "D48" "1.928"
$ws.Range("E48").Value = "  -0.93%  "
Set-TextValue "D49" "0.06817"
$ws.Range("E49").Value = "  -1.56%  "
Set-TextValue "D50" "1.122"
Set-TextValue "D51" "72.05"
$ws.Range("E51").Value = "  -0.08%  "
